$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the confidential disclosure text (date 2021-05-24 -> 2021-05-25)
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-25 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-8
$ws.Range("D2").Value = 0.5001774773288865
$ws.Range("E2").Value = -0.005616474993313747

$ws.Range("D3").Value = 0.2456704272605177
$ws.Range("E3").Value = 0.001308710193398399

$ws.Range("D4").Value = 0.09494334026774374
$ws.Range("E4").Value = -0.002269575085109055

$ws.Range("D5").Value = 0.102515009565534
$ws.Range("E5").Value = -0.0115848007414272

$ws.Range("D6").Value = 0.0298994193793082
$ws.Range("E6").Value = -0.0158245948522403

$ws.Range("D7").Value = 0.02679432619800997
$ws.Range("E7").Value = -0.005616662766206315

$ws.Range("E8").Value = -0.004514460791951613

# Restore sheet protection (content was protected before this edit)
$ws.Protect()
